$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-match data (country/tournament/season/date in A:E stay put).
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-RowData {
    param($ws, $r1, $r2)
    $vals1 = @{}
    $vals2 = @{}
    foreach ($c in $cols) {
        $vals1[$c] = $ws.Range("$c$r1").Value2
        $vals2[$c] = $ws.Range("$c$r2").Value2
    }
    foreach ($c in $cols) {
        $ws.Range("$c$r1").Value = $vals2[$c]
        $ws.Range("$c$r2").Value = $vals1[$c]
    }
}

# The rows below had their match data (F:V) swapped between the two
# neighbouring rows (row order of two fixtures on the same matchday flipped).
Swap-RowData $ws 72 73
Swap-RowData $ws 118 119
Swap-RowData $ws 120 121
Swap-RowData $ws 130 131
Swap-RowData $ws 147 148

# A new fixture row (Ferrol vs Andorra) was appended at the end of the table.
$ws.Range("A149:V149").Copy()
$ws.Range("A150:V150").PasteSpecial(-4122)

$ws.Range("A150").Value = 149
$ws.Range("B150").Value = "spain"
$ws.Range("C150").Value = "laliga2"
$ws.Range("D150").Value = "2023-2024"
$ws.Range("E150").Value = 45235.58333333334
$ws.Range("F150").Value = "Ferrol"
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = "Andorra"
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 2.04
$ws.Range("K150").Value = "29/10/2023 21:12"
$ws.Range("L150").Value = 2.23
$ws.Range("M150").Value = "05/11/2023 13:53"
$ws.Range("N150").Value = 3.3
$ws.Range("O150").Value = "29/10/2023 21:12"
$ws.Range("P150").Value = 3.15
$ws.Range("Q150").Value = "05/11/2023 13:53"
$ws.Range("R150").Value = 4.04
$ws.Range("S150").Value = "29/10/2023 21:12"
$ws.Range("T150").Value = 3.83
$ws.Range("U150").Value = "05/11/2023 13:53"
$ws.Range("V150").Value = "https://www.betexplorer.com/football/spain/laliga2/ferrol-fc-andorra/AsONbMm2/"

Write-Output "edit complete"
